$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '91.279.74'
Set-TextValue 'E2' '  +0.95%  '

Set-TextValue 'D3' '3.153.90'
Set-TextValue 'E3' '  +1.77%  '

Set-TextValue 'E4' '  +0.06%  '

Set-TextValue 'D5' '242.89'
Set-TextValue 'E5' '  +2.12%  '

Set-TextValue 'D6' '617.99'
Set-TextValue 'E6' '  -1.01%  '

Set-TextValue 'D7' '1.13'
Set-TextValue 'E7' '  -0.43%  '

Set-TextValue 'E8' '  +1.28%  '

Set-TextValue 'E9' '  +0.00%  '

Set-TextValue 'D10' '3.153.53'
Set-TextValue 'E10' '  +1.72%  '

Set-TextValue 'D11' '0.740'
Set-TextValue 'E11' '  +0.45%  '

Set-TextValue 'D12' '0.204'
Set-TextValue 'E12' '  +0.95%  '

Set-TextValue 'E13' '  -0.34%  '

Set-TextValue 'D14' '5.66'
Set-TextValue 'E14' '  +4.01%  '

Set-TextValue 'D15' '35.20'
Set-TextValue 'E15' '  +0.48%  '

Set-TextValue 'D16' '90.895.43'
Set-TextValue 'E16' '  +0.59%  '

Set-TextValue 'D17' '3.741.06'
Set-TextValue 'E17' '  +1.30%  '

Set-TextValue 'D18' '3.165.62'
Set-TextValue 'E18' '  +1.63%  '

Set-TextValue 'D19' '3.72'
Set-TextValue 'E19' '  -2.72%  '

Set-TextValue 'D20' '15.14'
Set-TextValue 'E20' '  +6.36%  '

Set-TextValue 'D21' '5.96'
Set-TextValue 'E21' '  +4.35%  '

Set-TextValue 'D22' '457.77'
Set-TextValue 'E22' '  +2.79%  '

Set-TextValue 'D23' '0.0000206'
Set-TextValue 'E23' '  -2.47%  '

Set-TextValue 'E24' '  +1.84%  '

Set-TextValue 'D25' '5.96'
Set-TextValue 'E25' '  +1.50%  '

Set-TextValue 'D26' '89.03'
Set-TextValue 'E26' '  -1.48%  '

Set-TextValue 'E27' '  -0.42%  '

Set-TextValue 'B28' 'WrappedeETH'
Set-TextValue 'C28' 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue 'D28' '3.327.16'
Set-TextValue 'E28' '  +2.03%  '

Set-TextValue 'B29' 'Hedera'
Set-TextValue 'C29' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D29' '0.146'
Set-TextValue 'E29' '  +34.89%  '

Set-TextValue 'E30' '  -0.03%  '

Set-TextValue 'D31' '0.235'
Set-TextValue 'E31' '  +5.57%  '

Set-TextValue 'D32' '0.169'
Set-TextValue 'E32' '  -4.52%  '

Set-TextValue 'D33' '9.41'
Set-TextValue 'E33' '  +2.73%  '

Set-TextValue 'D34' '0.174'
Set-TextValue 'E34' '  +12.83%  '

Set-TextValue 'B35' 'EthereumClassic'
Set-TextValue 'C35' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D35' '26.57'
Set-TextValue 'E35' '  +0.82%  '

Set-TextValue 'B36' 'RenderToken'
Set-TextValue 'C36' 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue 'D36' '7.61'
Set-TextValue 'E36' '  +3.22%  '

Set-TextValue 'B37' 'PancakeSwap'
Set-TextValue 'C37' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D37' '1.94'
Set-TextValue 'E37' '  +1.15%  '

Set-TextValue 'B38' 'Bittensor'
Set-TextValue 'C38' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D38' '498.84'
Set-TextValue 'E38' '  +1.30%  '

Set-TextValue 'B39' 'Fetch.AI'
Set-TextValue 'C39' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D39' '1.33'
Set-TextValue 'E39' '  +3.72%  '

Set-TextValue 'B40' 'MantraDAO'
Set-TextValue 'C40' 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextValue 'D40' '3.86'
Set-TextValue 'E40' '  -8.23%  '

Set-TextValue 'B41' 'PolygonEcosystemToken'
Set-TextValue 'C41' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D41' '0.447'
Set-TextValue 'E41' '  +7.36%  '

Set-TextValue 'B42' 'dogwifhat'
Set-TextValue 'C42' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D42' '3.49'
Set-TextValue 'E42' '  -2.60%  '

Set-TextValue 'B43' 'WhiteBITCoin'
Set-TextValue 'C43' 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 'D43' '22.13'
Set-TextValue 'E43' '  +0.15%  '

Set-TextValue 'B44' 'USDe'
Set-TextValue 'C44' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D44' '1.00'
Set-TextValue 'E44' '  +0.01%  '

Set-TextValue 'B45' 'Binance-PegBSC-USD'
Set-TextValue 'C45' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D45' '0.703'
Set-TextValue 'E45' '  -29.53%  '

Set-TextValue 'D46' '0.714'
Set-TextValue 'E46' '  +5.56%  '

Set-TextValue 'D47' '1.93'
Set-TextValue 'E47' '  +2.00%  '

Set-TextValue 'D48' '155.34'
Set-TextValue 'E48' '  -1.28%  '

Set-TextValue 'D49' '1.38'
Set-TextValue 'E49' '  +4.02%  '

Set-TextValue 'D50' '4.51'
Set-TextValue 'E50' '  -0.58%  '

Set-TextValue 'B51' 'VeChain'
Set-TextValue 'C51' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D51' '0.0328'
Set-TextValue 'E51' '  +8.83%  '
